# Generate Report for Handoff
# Updates the localization-status report to reflect that the content is
# now "Ready for handoff" instead of "In Translation", and refreshes the
# handoff timestamps accordingly.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns for zh-cn (E2) and de-de (F2)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-21 22:49:19"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-21 22:49:24"

# The "Status" columns widen to fit the new, longer text (target stored
# OOXML column width ~17.216; expressed here in COM "characters" units,
# which this host rounds to the nearest 1/6th character).
$newStatusWidth = 16.333333333333332
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth
